# Datenbankstruktur.xlsx update:
#  - mark every table-header band with a small "x" flag cell in column I
#    (used elsewhere in the workbook as a quick visual marker)
#  - rename the two foreign-key columns of the new "Spielplan" junction
#    table (row 51) from the generic "Verein_ID" to the explicit
#    "Heim_Verein_ID" / "Gast_Verein_ID" so home/away club references are
#    distinguishable
#  - scroll the sheet down a little and leave the selection on the newly
#    edited cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that carry a merged "table title" band (A<row>:H<row>). Each one
# gets a new flag cell in column I with the value "x".
$headerRows = @(2, 6, 10, 14, 18, 22, 26, 30, 34, 38, 42, 46)
foreach ($r in $headerRows) {
    $ws.Cells.Item($r, 9).Value = "x"
}

# New dummy-table columns: home club id / away club id.
$ws.Range("D51").Value = "Heim_Verein_ID"
$ws.Range("E51").Value = "Gast_Verein_ID"

# Update the view: scroll a little further down and select the cell that
# was just edited.
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E51").Select()
